$wb = $excel.ActiveWorkbook

$ov = $wb.Worksheets.Item("Overview")
$zh = $wb.Worksheets.Item("zh-cn")
$de = $wb.Worksheets.Item("de-de")

# Generate Report for Handoff:
# Rows 4-7 (0b787158-..., 0cef30ca-..., 44681964-..., 9ade4afb-...) move from
# "low" priority / "Ready for handoff" into the handed-off state: priority
# becomes "ht" and the Latest Handoff Datetime is refreshed to the new
# handoff timestamp, for both the zh-cn and de-de locale sheets. The de-de
# handoff timestamp is shared with the Overview sheet's "Latest HO Xliff
# Generate Date" column, so it is refreshed there too.

$zhHandoffTime = "2016-08-23 16:34:18"
$deHandoffTime = "2016-08-23 16:34:23"

foreach ($row in 4..7) {
    $zh.Range("E$row").Value = "ht"
    $zh.Range("H$row").Value = $zhHandoffTime

    $de.Range("E$row").Value = "ht"
    $de.Range("H$row").Value = $deHandoffTime

    $ov.Range("G$row").Value = $deHandoffTime
}
